$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds data in rows 2..53 (row 1 is the header).
# Two new product rows are being inserted at rows 8 and 9, pushing the
# existing rows 8..53 down to rows 10..55.
#
# We avoid Range.EntireRow.Insert() because (in this engine) it leaves a
# spare, unused cell-style entry behind. Instead we shift the data
# manually (bottom-up, so we never clobber a row before reading it) using
# Value2 (keeps numbers/booleans/strings/shared-strings intact), then
# write the two brand-new rows on top of the now-vacated rows 8 and 9.

for ($r = 53; $r -ge 8; $r--) {
    $src = $ws.Range("A" + $r + ":P" + $r)
    $vals = $src.Value2
    $dst = $ws.Range("A" + ($r + 2) + ":P" + ($r + 2))
    $dst.Value2 = $vals
}

# Rows 54 and 55 are brand new (the sheet previously ended at row 53), so
# column A's "barcode" number format (style index 1 in the original file)
# has to be (re)applied explicitly - it isn't inherited automatically the
# way column O's column-level style is.
$ws.Range("A54").NumberFormat = "0"
$ws.Range("A55").NumberFormat = "0"

# New row 8: Crema dental Colgate anticaries herbal
$ws.Range("A8").Value2 = 7509546695518
$ws.Range("B8").Value2 = "Crema dental"
$ws.Range("C8").Value2 = "anticaries"
$ws.Range("D8").Value2 = "herbal"
$ws.Range("E8").Value2 = "Colgate"
$ws.Range("F8").Value2 = 90
$ws.Range("G8").Value2 = "gr."
$ws.Range("H8").Value2 = "caja"
$ws.Range("I8").Value2 = "Cremas dentales"
$ws.Range("J8").Value2 = "Argentina"
$ws.Range("K8").Value2 = 6
$ws.Range("L8").Value2 = $false
$ws.Range("M8").Value2 = $true
$ws.Range("N8").Value2 = "C:\VentaSoft\Imágenes de artículos\7509546695518.png"
$ws.Range("O8").Value2 = $true
$ws.Range("P8").Value2 = $false

# New row 9: Azúcar común tipo "a" Bella Vista
$ws.Range("A9").Value2 = 7790220000746
$ws.Range("B9").Value2 = "Azúcar"
$ws.Range("C9").Value2 = "común"
$ws.Range("D9").Value2 = "tipo ""a"""
$ws.Range("E9").Value2 = "Bella Vista"
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = "kg."
$ws.Range("H9").Value2 = "bolsa"
$ws.Range("I9").Value2 = "Azúcar"
$ws.Range("J9").Value2 = "Argentina"
$ws.Range("K9").Value2 = 10
$ws.Range("L9").Value2 = $false
$ws.Range("M9").Value2 = $true
$ws.Range("N9").Value2 = "C:\VentaSoft\Imágenes de artículos\7790220000746.png"
$ws.Range("O9").Value2 = $true
$ws.Range("P9").Value2 = $true
